$wb = $excel.ActiveWorkbook

# --- Sheet "eps_factors" (first sheet) ---
$wsEps = $wb.Worksheets.Item("eps_factors")
$wsEps.Range("A24").Value = 0.119521198487387
$wsEps.Range("B24").Value = 0.120600956941548
$wsEps.Activate()
$excel.Application.Goto($wsEps.Range("B2:B26"))

# --- Sheet "s_factors" (second sheet) ---
$wsS = $wb.Worksheets.Item("s_factors")
$wsS.Range("A24").Value = 0.404433
$wsS.Range("B24").Value = 0.2579358
$wsS.Activate()
$excel.Application.Goto($wsS.Range("B2:B26"))

# Make "s_factors" the active sheet/tab, matching activeTab="1" in workbook.xml
$wsS.Activate()
